$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (824cce6b...) status text changes ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Handback transform failed"
$ov.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: widen "Error Detail" column (P), refresh Status (C), and set error detail on row 3 ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Columns.Item(16).ColumnWidth = 39.166666666666664
$zh.Range("C3").Value = "Handback transform failed"
$zh.Range("P3").Value = "Handback file name: f4ttkhot.tbl is different with handoff file name: 824cce6b-40bc-44d0-b46b-aaaf920e5a05.64cafc2d706588a02179ad98f55ec3cef40c5b96.zh-cn."

# --- de-de sheet: widen "Error Detail" column (P), refresh Status (C), and set error detail on row 3 ---
$de = $wb.Worksheets.Item("de-de")
$de.Columns.Item(16).ColumnWidth = 39.166666666666664
$de.Range("C3").Value = "Handback transform failed"
$de.Range("P3").Value = "Handback file name: f4ttkhot.tbl is different with handoff file name: 824cce6b-40bc-44d0-b46b-aaaf920e5a05.64cafc2d706588a02179ad98f55ec3cef40c5b96.de-de."
